$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes: "% Avance" progress updated, note added about the missing account ---

# G8: progress went from 100% to 70%
$ws.Range("G8").Value = 0.7

# H8: new note added next to row 8 explaining what's still missing
$ws.Range("H8").Value = "Falta Crear cuenta"

# G9: progress went from 40% to 80%
$ws.Range("G9").Value = 0.8

# --- View state: scrolled/selected a different cell while reviewing the sheet ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F10").Select()
